$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.925179770859371
$ws.Range("C2").Value = 5.780391821375774
$ws.Range("D2").Value = 4.443596211984586
$ws.Range("F2").Value = 19.58795638235258
$ws.Range("G2").Value = 20.57547636414088
$ws.Range("H2").Value = 12.38124052643509
$ws.Range("I2").Value = 17.45920573582139
$ws.Range("K2").Value = 8.471841746582419
$ws.Range("N2").Value = 16.86580294599075
$ws.Range("O2").Value = 17.72951726136178
$ws.Range("B3").Value = 7.577721788249709
$ws.Range("C3").Value = 5.653444319448128
$ws.Range("D3").Value = 4.324653435168544
$ws.Range("F3").Value = 19.64220381829719
$ws.Range("G3").Value = 20.67328514877434
$ws.Range("H3").Value = 12.42391561314609
$ws.Range("I3").Value = 17.54402064359958
$ws.Range("K3").Value = 8.205932179077227
$ws.Range("N3").Value = 16.90251082953466
$ws.Range("O3").Value = 17.80530891082829
$ws.Range("B4").Value = 7.356420937903089
$ws.Range("C4").Value = 5.573449165194751
$ws.Range("D4").Value = 4.249076428820368
$ws.Range("F4").Value = 19.6806658384917
$ws.Range("G4").Value = 20.74036784554328
$ws.Range("H4").Value = 12.45181086161699
$ws.Range("I4").Value = 17.59913069306899
$ws.Range("K4").Value = 8.036753638461693
$ws.Range("N4").Value = 16.92649702482272
$ws.Range("O4").Value = 17.85532442378621
$ws.Range("B5").Value = 7.264357016370989
$ws.Range("C5").Value = 5.540365615358892
$ws.Range("D5").Value = 4.217663963746512
$ws.Range("F5").Value = 19.69763251308434
$ws.Range("G5").Value = 20.76946356024395
$ws.Range("H5").Value = 12.46360451203502
$ws.Range("I5").Value = 17.62235220767446
$ws.Range("K5").Value = 7.966392283048599
$ws.Range("N5").Value = 16.93663625174448
$ws.Range("O5").Value = 17.87658056519112
$ws.Range("B6").Value = 7.248960026703576
$ws.Range("C6").Value = 5.534843750826802
$ws.Range("D6").Value = 4.212411648058549
$ws.Range("F6").Value = 19.70052782530942
$ws.Range("G6").Value = 20.77440087417853
$ws.Range("H6").Value = 12.46558858954096
$ws.Range("I6").Value = 17.62625427810901
$ws.Range("K6").Value = 7.954625144005741
$ws.Range("N6").Value = 16.93834190765768
$ws.Range("O6").Value = 17.88016293108266
$ws.Range("B7").Value = 7.355186781622431
$ws.Range("C7").Value = 5.573004910817899
$ws.Range("D7").Value = 4.248655240552278
$ws.Range("F7").Value = 19.68088942536939
$ws.Range("G7").Value = 20.74075312953896
$ws.Range("H7").Value = 12.45196818893488
$ws.Range("I7").Value = 17.5994407729539
$ws.Range("K7").Value = 8.035810376965575
$ws.Range("N7").Value = 16.92663228856648
$ws.Range("O7").Value = 17.85560755164114
$ws.Range("B8").Value = 7.8070951573112
$ws.Range("C8").Value = 5.737063240161848
$ws.Range("D8").Value = 4.4031279939298
$ws.Range("F8").Value = 19.60558956129119
$ws.Range("G8").Value = 20.6077373810518
$ws.Range("H8").Value = 12.39560390291545
$ws.Range("I8").Value = 17.48782092060422
$ws.Range("K8").Value = 8.381420250194797
$ws.Range("N8").Value = 16.87815983253931
$ws.Range("O8").Value = 17.75492786185322
$ws.Range("B9").Value = 8.625665451152484
$ws.Range("C9").Value = 6.041279496743508
$ws.Range("D9").Value = 4.684817473937795
$ws.Range("F9").Value = 19.49893733558453
$ws.Range("G9").Value = 20.40300991427211
$ws.Range("H9").Value = 12.29848085032435
$ws.Range("I9").Value = 17.29295896800454
$ws.Range("K9").Value = 9.009689021154019
$ws.Range("N9").Value = 16.79455876525473
$ws.Range("O9").Value = 17.58512028192538
$ws.Range("B10").Value = 9.180863060964782
$ws.Range("C10").Value = 6.25257714085989
$ws.Range("D10").Value = 4.877592974897566
$ws.Range("F10").Value = 19.44572182454
$ws.Range("G10").Value = 20.28727132525483
$ws.Range("H10").Value = 12.23526633895243
$ws.Range("I10").Value = 17.1643778042528
$ws.Range("K10").Value = 9.438080004748912
$ws.Range("N10").Value = 16.74007623699762
$ws.Range("O10").Value = 17.47723018376691
$ws.Range("B11").Value = 9.422596955943611
$ws.Range("C11").Value = 6.345760641056568
$ws.Range("D11").Value = 4.961988905182224
$ws.Range("F11").Value = 19.42699485782058
$ws.Range("G11").Value = 20.24224058338872
$ws.Range("H11").Value = 12.20826980949962
$ws.Range("I11").Value = 17.10903698996748
$ws.Range("K11").Value = 9.625226660831727
$ws.Range("N11").Value = 16.7167885128073
$ws.Range("O11").Value = 17.43181685799047
$ws.Range("B12").Value = 9.512524568908999
$ws.Range("C12").Value = 6.380602639091205
$ws.Range("D12").Value = 4.993456838343261
$ws.Range("F12").Value = 19.42069283320215
$ws.Range("G12").Value = 20.22629127962427
$ws.Range("H12").Value = 12.19829956980279
$ws.Range("I12").Value = 17.08853319409509
$ws.Range("K12").Value = 9.694944921778491
$ws.Range("N12").Value = 16.70818460889065
$ws.Range("O12").Value = 17.41514784834153
$ws.Range("B13").Value = 9.493229510054141
$ws.Range("C13").Value = 6.373118899433448
$ws.Range("D13").Value = 4.986701739000083
$ws.Range("F13").Value = 19.4220149586309
$ws.Range("G13").Value = 20.22967708739808
$ws.Range("H13").Value = 12.20043560295971
$ws.Range("I13").Value = 17.09292893440688
$ws.Range("K13").Value = 9.679981562464771
$ws.Range("N13").Value = 16.71002807572301
$ws.Range("O13").Value = 17.41871432056032
$ws.Range("B14").Value = 9.430027939030508
$ws.Range("C14").Value = 6.348636144271155
$ws.Range("D14").Value = 4.96458772009764
$ws.Range("F14").Value = 19.42646055827181
$ws.Range("G14").Value = 20.24090628185083
$ws.Range("H14").Value = 12.20744448807112
$ws.Range("I14").Value = 17.10734106347749
$ws.Range("K14").Value = 9.63098567820221
$ws.Range("N14").Value = 16.71607636570746
$ws.Range("O14").Value = 17.43043490135415
$ws.Range("B15").Value = 9.391103779863236
$ws.Range("C15").Value = 6.333581223190541
$ws.Range("D15").Value = 4.95097784094301
$ws.Range("F15").Value = 19.42928646058457
$ws.Range("G15").Value = 20.24792831927191
$ws.Range("H15").Value = 12.21177053964566
$ws.Range("I15").Value = 17.11622782941692
$ws.Range("K15").Value = 9.600823442126396
$ws.Range("N15").Value = 16.71980905539343
$ws.Range("O15").Value = 17.43768288483606
$ws.Range("B16").Value = 9.16484241183629
$ws.Range("C16").Value = 6.246426340050403
$ws.Range("D16").Value = 4.872009664494907
$ws.Range("F16").Value = 19.44705608922317
$ws.Range("G16").Value = 20.29036813678229
$ws.Range("H16").Value = 12.2370660137657
$ws.Range("I16").Value = 17.1680578189379
$ws.Range("K16").Value = 9.42569057721636
$ws.Range("N16").Value = 16.74162821245749
$ws.Range("O16").Value = 17.48027189766585
$ws.Range("B17").Value = 9.023223228350519
$ws.Range("C17").Value = 6.192191044424479
$ws.Range("D17").Value = 4.822708007057601
$ws.Range("F17").Value = 19.45936206374698
$ws.Range("G17").Value = 20.31836034946156
$ws.Range("H17").Value = 12.2530345211054
$ws.Range("I17").Value = 17.20066055160516
$ws.Range("K17").Value = 9.316243394330561
$ws.Range("N17").Value = 16.75539648149916
$ws.Range("O17").Value = 17.50733850451699
$ws.Range("B18").Value = 8.940751066146028
$ws.Range("C18").Value = 6.160721122429552
$ws.Range("D18").Value = 4.794041398750956
$ws.Range("F18").Value = 19.46695598344815
$ws.Range("G18").Value = 20.3351774918209
$ws.Range("H18").Value = 12.26238486951252
$ws.Range("I18").Value = 17.21970943960008
$ws.Range("K18").Value = 9.252566623794452
$ws.Range("N18").Value = 16.7634565306568
$ws.Range("O18").Value = 17.52325155293876
$ws.Range("B19").Value = 8.912654628491683
$ws.Range("C19").Value = 6.150019386389004
$ws.Range("D19").Value = 4.784282738953192
$ws.Range("F19").Value = 19.4696157052156
$ws.Range("G19").Value = 20.3409943650062
$ws.Range("H19").Value = 12.26557920385862
$ws.Range("I19").Value = 17.22621003595044
$ws.Range("K19").Value = 9.230883381894527
$ws.Range("N19").Value = 16.76620974276416
$ws.Range("O19").Value = 17.52869867641556
$ws.Range("B20").Value = 9.038404465860211
$ws.Range("C20").Value = 6.197993134706245
$ws.Range("D20").Value = 4.827988420913701
$ws.Range("F20").Value = 19.45799867260799
$ws.Range("G20").Value = 20.31530629484912
$ws.Range("H20").Value = 12.25131750089847
$ws.Range("I20").Value = 17.19715923646575
$ws.Range("K20").Value = 9.327969628763048
$ws.Range("N20").Value = 16.75391624674802
$ws.Range("O20").Value = 17.50442150063409
$ws.Range("B21").Value = 9.448635907657996
$ws.Range("C21").Value = 6.355839557074441
$ws.Range("D21").Value = 4.971096592570325
$ws.Range("F21").Value = 19.42513334351521
$ws.Range("G21").Value = 20.23757800862395
$ws.Range("H21").Value = 12.20537895139613
$ws.Range("I21").Value = 17.10309559380503
$ws.Range("K21").Value = 9.645408454743555
$ws.Range("N21").Value = 16.71429401518428
$ws.Range("O21").Value = 17.42697794558183
$ws.Range("B22").Value = 9.707331845569863
$ws.Range("C22").Value = 6.45640092929954
$ws.Range("D22").Value = 5.061757741909593
$ws.Range("F22").Value = 19.40825590566094
$ws.Range("G22").Value = 20.19320942598469
$ws.Range("H22").Value = 12.17682858559439
$ws.Range("I22").Value = 17.0442572340542
$ws.Range("K22").Value = 9.846154637983316
$ws.Range("N22").Value = 16.68964949716409
$ws.Range("O22").Value = 17.37944225664295
$ws.Range("B23").Value = 9.570138188409134
$ws.Range("C23").Value = 6.40297426679811
$ws.Range("D23").Value = 5.01363769378727
$ws.Range("F23").Value = 19.41684230861556
$ws.Range("G23").Value = 20.21629903761312
$ws.Range("H23").Value = 12.19193176931672
$ws.Range("I23").Value = 17.0754192157368
$ws.Range("K23").Value = 9.739638900021234
$ws.Range("N23").Value = 16.702688466495
$ws.Range("O23").Value = 17.40453102874744
$ws.Range("B24").Value = 9.031544305269195
$ws.Range("C24").Value = 6.195370909187081
$ws.Range("D24").Value = 4.82560214882774
$ws.Range("F24").Value = 19.45861344506853
$ws.Range("G24").Value = 20.31668477817413
$ws.Range("H24").Value = 12.25209323667173
$ws.Range("I24").Value = 17.19874123050838
$ws.Range("K24").Value = 9.322670544162975
$ws.Range("N24").Value = 16.75458501077099
$ws.Range("O24").Value = 17.50573918141041
$ws.Range("B25").Value = 8.412024277236705
$ws.Range("C25").Value = 5.961020664294655
$ws.Range("D25").Value = 4.611021469501335
$ws.Range("F25").Value = 19.52338318624341
$ws.Range("G25").Value = 20.45233719800569
$ws.Range("H25").Value = 12.32332318958954
$ws.Range("I25").Value = 17.34310885345517
$ws.Range("K25").Value = 9.009689021154019
$ws.Range("N25").Value = 16.81595340958717
$ws.Range("O25").Value = 17.62809759770606
